$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 holds the single data record on this sheet.
# Update the Client Id, Candidate ID, User Name, Exam Password,
# First Name and Last Name fields as per the "Added Modified Reg
# iExam P2,P3 TC's" commit.

$ws.Range("A2").Value = "TlcTL154"
$ws.Range("B2").Value = 23080719
$ws.Range("C2").Value = "itfhqke80"
$ws.Range("D2").Value = "q8&eG3#E"
$ws.Range("F2").Value = "AgynuRYv"
$ws.Range("G2").Value = "XvtW"
